$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.355.84'
$ws.Range('E2').Value = '  -0.12%  '

$ws.Range('D3').Value = '1.627.05'
$ws.Range('E3').Value = '  -0.03%  '

$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.39%  '

$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.46%  '

$ws.Range('D6').Value = '303.64'
$ws.Range('E6').Value = '  -1.09%  '

$ws.Range('D7').Value = '0.3771'
$ws.Range('E7').Value = '  -0.33%  '

$ws.Range('D8').Value = '52.10'
$ws.Range('E8').Value = '  -2.47%  '

$ws.Range('D9').Value = '0.3621'
$ws.Range('E9').Value = '  -1.12%  '

$ws.Range('D10').Value = '1.236'
$ws.Range('E10').Value = '  -3.28%  '

$ws.Range('D11').Value = '0.08087'
$ws.Range('E11').Value = '  -1.26%  '

$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.38%  '

$ws.Range('D13').Value = '22.69'
$ws.Range('E13').Value = '  -2.20%  '

$ws.Range('D14').Value = '6.565'
$ws.Range('E14').Value = '  -1.38%  '

$ws.Range('E15').Value = '  -1.14%  '

$ws.Range('D16').Value = '7.228'
$ws.Range('E16').Value = '  -3.04%  '

$ws.Range('D17').Value = '1.627.53'

$ws.Range('D18').Value = '93.48'
$ws.Range('E18').Value = '  -1.31%  '

$ws.Range('D19').Value = '0.06925'
$ws.Range('E19').Value = '  -0.03%  '

$ws.Range('D20').Value = '17.95'
$ws.Range('E20').Value = '  -2.24%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.48%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.458'
$ws.Range('E22').Value = '  -1.86%  '

$ws.Range('D23').Value = '23.354.08'
$ws.Range('E23').Value = '  -0.22%  '

$ws.Range('D24').Value = '12.73'
$ws.Range('E24').Value = '  -1.90%  '

$ws.Range('D25').Value = '3.223'
$ws.Range('E25').Value = '  +2.58%  '

$ws.Range('D26').Value = '2.446'
$ws.Range('E26').Value = '  +1.18%  '

$ws.Range('D27').Value = '21.07'
$ws.Range('E27').Value = '  -1.41%  '

$ws.Range('D28').Value = '149.91'
$ws.Range('E28').Value = '  -0.47%  '

$ws.Range('D29').Value = '5.296'
$ws.Range('E29').Value = '  +0.43%  '

$ws.Range('D30').Value = '134.87'
$ws.Range('E30').Value = '  -1.13%  '

$ws.Range('D31').Value = '2.301'
$ws.Range('E31').Value = '  -4.71%  '

$ws.Range('D32').Value = '1.808.92'
$ws.Range('E32').Value = '  +0.43%  '

$ws.Range('D33').Value = '6.772'
$ws.Range('E33').Value = '  -1.88%  '

$ws.Range('E34').Value = '  +4.33%  '

$ws.Range('D35').Value = '0.9477'
$ws.Range('E35').Value = '  -2.82%  '

$ws.Range('D36').Value = '0.02823'
$ws.Range('E36').Value = '  +1.09%  '

$ws.Range('D37').Value = '0.2531'
$ws.Range('E37').Value = '  -0.06%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '6.123'
$ws.Range('E38').Value = '  -1.76%  '

$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '0.08826'
$ws.Range('E39').Value = '  -0.18%  '

$ws.Range('D40').Value = '0.07109'
$ws.Range('E40').Value = '  -4.52%  '

$ws.Range('D41').Value = '1.361'
$ws.Range('E41').Value = '  -3.23%  '

$ws.Range('D42').Value = '0.7046'
$ws.Range('E42').Value = '  -1.63%  '

$ws.Range('D43').Value = '16.20'
$ws.Range('E43').Value = '  +0.26%  '

$ws.Range('D44').Value = '12.35'
$ws.Range('E44').Value = '  -3.39%  '

$ws.Range('D45').Value = '0.6453'
$ws.Range('E45').Value = '  -2.43%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '2.319'
$ws.Range('E46').Value = '  -1.75%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('E48').Value = '  -1.11%  '

$ws.Range('D49').Value = '0.07978'
$ws.Range('E49').Value = '  -0.46%  '

$ws.Range('E50').Value = '  -0.77%  '

$ws.Range('D51').Value = '126.35'
$ws.Range('E51').Value = '  -3.77%  '
